$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-04-14 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-04-15 Monday", 2) | Out-Null
$d.Content.Find.Execute("23×40=", $true, $false, $false, $false, $false, $true, 1, $false, "75×97=", 2) | Out-Null
$d.Content.Find.Execute("30×33=", $true, $false, $false, $false, $false, $true, 1, $false, "93×79=", 2) | Out-Null
$d.Content.Find.Execute("98×83=", $true, $false, $false, $false, $false, $true, 1, $false, "61×80=", 2) | Out-Null
$d.Content.Find.Execute("22×28=", $true, $false, $false, $false, $false, $true, 1, $false, "38×38=", 2) | Out-Null
$d.Content.Find.Execute("91×91=", $true, $false, $false, $false, $false, $true, 1, $false, "63×41=", 2) | Out-Null
$d.Content.Find.Execute("14×72=", $true, $false, $false, $false, $false, $true, 1, $false, "94×92=", 2) | Out-Null
$d.Content.Find.Execute("38×62=", $true, $false, $false, $false, $false, $true, 1, $false, "84×69=", 2) | Out-Null
$d.Content.Find.Execute("76×65=", $true, $false, $false, $false, $false, $true, 1, $false, "99×58=", 2) | Out-Null
$d.Content.Find.Execute("87×19=", $true, $false, $false, $false, $false, $true, 1, $false, "68×46=", 2) | Out-Null
$d.Content.Find.Execute("49×47=", $true, $false, $false, $false, $false, $true, 1, $false, "62×19=", 2) | Out-Null
$d.Content.Find.Execute("72×88=", $true, $false, $false, $false, $false, $true, 1, $false, "85×34=", 2) | Out-Null
$d.Content.Find.Execute("40×20=", $true, $false, $false, $false, $false, $true, 1, $false, "27×75=", 2) | Out-Null
$d.Content.Find.Execute("26×30=", $true, $false, $false, $false, $false, $true, 1, $false, "43×90=", 2) | Out-Null
$d.Content.Find.Execute("48×38=", $true, $false, $false, $false, $false, $true, 1, $false, "76×98=", 2) | Out-Null
$d.Content.Find.Execute("15×14=", $true, $false, $false, $false, $false, $true, 1, $false, "31×75=", 2) | Out-Null
$d.Content.Find.Execute("17×82=", $true, $false, $false, $false, $false, $true, 1, $false, "96×73=", 2) | Out-Null
$d.Content.Find.Execute("89×64=", $true, $false, $false, $false, $false, $true, 1, $false, "68×30=", 2) | Out-Null
$d.Content.Find.Execute("18×43=", $true, $false, $false, $false, $false, $true, 1, $false, "34×41=", 2) | Out-Null
$d.Content.Find.Execute("41×13=", $true, $false, $false, $false, $false, $true, 1, $false, "86×15=", 2) | Out-Null
$d.Content.Find.Execute("40×90=", $true, $false, $false, $false, $false, $true, 1, $false, "17×56=", 2) | Out-Null
$d.Content.Find.Execute("97×69=", $true, $false, $false, $false, $false, $true, 1, $false, "12×17=", 2) | Out-Null
$d.Content.Find.Execute("73×22=", $true, $false, $false, $false, $false, $true, 1, $false, "73×65=", 2) | Out-Null
$d.Content.Find.Execute("43×81=", $true, $false, $false, $false, $false, $true, 1, $false, "31×95=", 2) | Out-Null
$d.Content.Find.Execute("33×51=", $true, $false, $false, $false, $false, $true, 1, $false, "99×70=", 2) | Out-Null
$d.Content.Find.Execute("36×48=", $true, $false, $false, $false, $false, $true, 1, $false, "41×75=", 2) | Out-Null
